$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the split "REQ" + "Txx"/"Lxx" runs into single runs, per the diff.
#    Find/Replace across a run boundary collapses the match into one run.
# ---------------------------------------------------------------------------
function Merge-ReqCell($tableIndex, $row, $col, $newText) {
    $cell = $d.Tables.Item($tableIndex).Cell($row, $col)
    $cell.Range.Find.Execute($newText, $false, $false, $false, $false, $false, `
                              $true, 1, $false, $newText, 2)
}

# "Para tiquetes" table (table 1): REQT1, REQT2, REQT3
Merge-ReqCell 1 1 1 "REQT1"
Merge-ReqCell 1 2 1 "REQT2"
Merge-ReqCell 1 3 1 "REQT3"

# "Para lineas" table (table 2): REQL1..REQL4 (REQL5 handled specially below)
Merge-ReqCell 2 1 1 "REQL1"
Merge-ReqCell 2 2 1 "REQL2"
Merge-ReqCell 2 3 1 "REQL3"
Merge-ReqCell 2 4 1 "REQL4"

# ---------------------------------------------------------------------------
# 2) REQL5 also needs a <w:lastRenderedPageBreak/> moved in front of its
#    text (it used to sit on the "Otros requerimientos:" run instead).
#    Merge the runs first, then splice the break into that paragraph.
# ---------------------------------------------------------------------------
Merge-ReqCell 2 5 1 "REQL5"

$reql5Para = $d.Tables.Item(2).Cell(5, 1).Range.Paragraphs.Item(1).Range
$reql5Xml = '<?xml version="1.0" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" ' + `
  'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>REQL5</w:t></w:r></w:p></w:body>' + `
  '</w:document></pkg:xmlData></pkg:part></pkg:package>'
$reql5Para.InsertXML($reql5Xml)

# Remove the old <w:lastRenderedPageBreak/> from the "Otros requerimientos:" run.
$otrosRange = $d.Content
$otrosRange.Find.ClearFormatting()
$otrosRange.Find.Execute("Otros requerimientos:", $false, $false, $false, $false, `
                          $false, $true, 1, $false, "Otros requerimientos:", 2)

# ---------------------------------------------------------------------------
# 3) Append three new rows to the end of the "tiquetes" table (table 1):
#    REQT6 / REQT7 / REQT8. REQT8's paragraph also gains the relocated
#    "_GoBack" bookmark (it used to live alone in the document's final
#    paragraph).
# ---------------------------------------------------------------------------
$tTable = $d.Tables.Item(1)

$row6 = $tTable.Rows.Add()
$row6.Cells.Item(1).Range.Text = "REQT6"
$row6.Cells.Item(2).Range.Text = "Ver los tiquetes actuales por línea."

$row7 = $d.Tables.Item(1).Rows.Add()
$row7.Cells.Item(1).Range.Text = "REQT7"
$row7.Cells.Item(2).Range.Text = "Ver los tiquetes totales."

$row8 = $d.Tables.Item(1).Rows.Add()
$row8.Cells.Item(1).Range.Text = "REQT8"
$row8.Cells.Item(2).Range.Text = "Ver un tiquete especifico."

$req8Para = $d.Tables.Item(1).Cell(8, 2).Range.Paragraphs.Item(1).Range
$req8Xml = '<?xml version="1.0" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" ' + `
  'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p><w:r><w:t>Ver un tiquete especifico.</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body>' + `
  '</w:document></pkg:xmlData></pkg:part></pkg:package>'
$req8Para.InsertXML($req8Xml)

# ---------------------------------------------------------------------------
# 4) Remove the old "_GoBack" bookmark paragraph at the very end of the
#    document, leaving a single empty paragraph in its place.
# ---------------------------------------------------------------------------
$tailPara = $d.Range($d.Content.End - 1, $d.Content.End)
$emptyParaXml = '<?xml version="1.0" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" ' + `
  'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$tailPara.InsertXML($emptyParaXml)

# The insert above splits the bookmark-only paragraph into two empty
# paragraphs (the fresh one + the old, now bookmark-free, one). Collapse
# back down to a single trailing empty paragraph.
$endNow = $d.Content.End
$newLast = $d.Range($endNow - 1, $endNow)
$newSecondLast = $d.Range($endNow - 2, $endNow - 1)
$dupRange = $d.Range($newSecondLast.Start, $newLast.Start)
$dupRange.Delete()

Write-Host "Done"
